$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 ("Docentes responsaveis:" value row, containing the
# professor name in B13/C13) is a duplicate of content that needs to move
# elsewhere, so the whole row is removed and everything below shifts up.
$ws.Rows(13).Delete()

# Two rows lose their custom row height entirely in the target layout
# (no ht/customHeight attribute at all), while keeping their single label
# cell untouched.
$ws.Rows(11).AutoFit()
$ws.Rows(16).AutoFit()

# One row needs a taller custom height than it had before (60 -> 120).
$ws.Rows(14).RowHeight = 120

# Now patch the B/C column contents that were reshuffled onto different
# rows (labels in column A already line up correctly after the delete).
$ws.Range("B10").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C10").Value = "1285870 - Marcos Villela Barcza"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2020" must land as literal text (it already exists as a shared
# string elsewhere), not get reinterpreted as a date serial number, so copy
# it from the existing B8:C8 cells that already hold it as text.
$ws.Range("B8:C8").Copy($ws.Range("B15:C15"))

$ws.Range("B18").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C18").Value = "1285870 - Marcos Villela Barcza"

$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."

$ws.Range("B20").Value = "Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula."

$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
